$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.132.48"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "1.898.95"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9979"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5185"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3757"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07253"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.12"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.34%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9010"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08360"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.66"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.37%  "

$ws.Range("D14").Value = "1.888.63"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.282"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9995"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008628"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9986"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.25%  "

$ws.Range("D20").Value = "27.171.07"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.078"
$ws.Range("D21").ClearFormats()

$ws.Range("D22").Value = "2.141.97"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.62"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.428"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.328"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.12"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.747"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.19"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.91"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.818"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.888"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09255"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05067"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7958"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.240"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.420"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.957"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.596"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5653"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01991"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.074"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.016"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.570"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.52"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1516"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4840"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.14"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9967"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.629"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.64"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.20%  "
